$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price-column (D) cells hold plain numeric-looking strings; force text
# format before assignment so Excel does not auto-convert them to numbers
# (matches the workbook's existing inlineStr/text storage for this column).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.144.16"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.221.07"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "293.79"
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.79"
$ws.Range("E6").Value = "  +0.84%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "30.69"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "50.88"
$ws.Range("E11").Value = "  +7.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0782"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("E13").Value = "  +3.94%  "
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.564.50"
$ws.Range("E15").Value = "  +0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.82"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.216.50"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.736"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "40.067.64"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0890"
$ws.Range("E20").Value = "  +1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.28"
$ws.Range("E21").Value = "  -2.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.78"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.66"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.01"
$ws.Range("E24").Value = "  +0.43%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.47"
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("E27").Value = "  +0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.17"
$ws.Range("E28").Value = "  +3.58%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.34"
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.07"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.18"
$ws.Range("E31").Value = "  +4.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.88"
$ws.Range("E32").Value = "  +1.06%  "
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.97"
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.04"
$ws.Range("E35").Value = "  +8.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0715"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  -2.06%  "
$ws.Range("E38").Value = "  +1.78%  "
$ws.Range("E39").Value = "  +4.44%  "
$ws.Range("E40").Value = "  +1.70%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "15.67"
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.079.36"
$ws.Range("E43").Value = "  -2.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.68"
$ws.Range("E44").Value = "  +13.28%  "
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.00"
$ws.Range("E46").Value = "  +0.83%  "
$ws.Range("E47").Value = "  +4.58%  "
$ws.Range("E48").Value = "  -10.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.436.20"
$ws.Range("E49").Value = "  +0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.13"
$ws.Range("E50").Value = "  +5.23%  "
$ws.Range("E51").Value = "  +2.31%  "
